# Update Leave Card 1/5/2024 4:46 PM
# - Shift the PERIOD dates in the leave table (rows 65-94) forward so the
#   monthly entries land on month-end dates instead of the 1st, and the
#   tail of the table becomes a run of consecutive daily entries in Feb 2024.
# - Record an EARNED credit of 1.25 for the 12/31/2023 period (row 77),
#   which ripples into the BALANCE totals (E9/I9) via the table formulas.
# - Scroll the frozen/split view down and move the active selection to
#   reflect where the user was last working in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- PERIOD column (Table1[PERIOD], column A) date updates ---------------
$dateUpdates = @{
    "A65" = 44957
    "A66" = 44985
    "A67" = 45016
    "A69" = 45046
    "A70" = 45077
    "A71" = 45107
    "A72" = 45138
    "A73" = 45169
    "A74" = 45199
    "A75" = 45230
    "A76" = 45260
    "A77" = 45291
    "A79" = 45322
    "A81" = 45324
    "A82" = 45325
    "A83" = 45326
    "A84" = 45327
    "A85" = 45328
    "A86" = 45329
    "A87" = 45330
    "A88" = 45331
    "A89" = 45332
    "A90" = 45333
    "A91" = 45334
    "A92" = 45335
    "A93" = 45336
    "A94" = 45337
}

foreach ($addr in $dateUpdates.Keys) {
    $ws.Range($addr).Value2 = $dateUpdates[$addr]
}

# --- EARNED credit added for the 12/31/2023 row (row 77) -----------------
$ws.Range("C77").Value2 = 1.25

# --- Scroll the split view so row 64 is the top of the lower pane, and ---
# --- leave the active selection on B79, matching where editing resumed. --
$ws.Activate()
$win = $excel.ActiveWindow
$win.SplitRow = 63
$ws.Range("B79").Select()
